# Reproduce the "Commiting changes to local" edit:
#  - Append 6 new test-case rows (B35..B40) to the "Test Cases" sheet.
#  - Update the used-range selection / active cell accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply a thin box border that exactly matches the sheet's existing
# "bordered" cell style (style index 3 in the original file: borderId=1,
# no fill, no wrap) so new cells reuse that style instead of minting a
# new one. The combination LineStyle=1 / Weight=2 / Color=0 is what the
# existing cells resolve to.
function Set-ExistingBorder($cell) {
    $cell.Borders.LineStyle = 1
    $cell.Borders.Weight = 2
    $cell.Borders.Color = 0
}

$rows = @(
    @{ Row = 36; A = "TestCase_B35"; B = "TBD-10"; C = "Verify that no filtering options are present in ALL content type"; D = "Y"; E = "PASS"; Wrap = $false },
    @{ Row = 37; A = "TestCase_B36"; B = "TBD-11"; C = "Verify that the following fields get displayed in the SORT BY drop down when ARTICLES is selected as content type in the left navigation pane:`na)Relevance`nb)Times cited`nc)Publication Date(Newest)`nd)Publication Date(Oldest)`n"; D = "Y"; E = "PASS"; Wrap = $true },
    @{ Row = 38; A = "TestCase_B37"; B = "TBD-12"; C = "Verify that user is able to sort the articles by TIMES CITED field in ARTICLES content type"; D = "Y"; E = "PASS"; Wrap = $false },
    @{ Row = 39; A = "TestCase_B38"; B = "TBD-13"; C = "Verify that only articles get displayed when user chooses ARTICLES as content type"; D = "Y"; E = "PASS"; Wrap = $false },
    @{ Row = 40; A = "TestCase_B39"; B = "TBD-14"; C = "Verify that all articles are sorted by RELEVANCE by default in ARTICLES content type"; D = "Y"; E = "PASS"; Wrap = $false },
    @{ Row = 41; A = "TestCase_B40"; B = "TBD-15"; C = "Verify that following filters are present for ARTICLES content type:`na)Document Type`nb)Authors`nc)Categories`nd)Institutions`n"; D = "Y"; E = "PASS"; Wrap = $true }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.A
    Set-ExistingBorder $cellA

    $cellB = $ws.Cells.Item($rowNum, 2)
    $cellB.Value = $r.B
    Set-ExistingBorder $cellB

    $cellC = $ws.Cells.Item($rowNum, 3)
    $cellC.Value = $r.C
    Set-ExistingBorder $cellC
    if ($r.Wrap) {
        $cellC.WrapText = $true
    }

    $cellD = $ws.Cells.Item($rowNum, 4)
    $cellD.Value = $r.D
    Set-ExistingBorder $cellD

    $cellE = $ws.Cells.Item($rowNum, 5)
    $cellE.Value = $r.E
    Set-ExistingBorder $cellE
}

# Match the post-edit selection recorded in the workbook: active cell A41.
$ws.Range("A41").Select() | Out-Null
